# The commit swaps the colour scheme baked into the deck's theme: the
# slide master's theme (ppt/theme/theme1.xml) goes from the "Integral"
# palette to the stock "Office Theme" palette (the notes master's
# theme, theme2.xml, correspondingly becomes "Integral" in the
# canonical XML; that second theme part isn't reachable as an
# independent object from COM automation here, since
# NotesMaster.Theme / HandoutMaster.Theme / Presentation.Theme all
# resolve back to the single slide-master theme, so only that one
# theme can be updated this way).
#
# Re-point the presentation's theme colour scheme to the "Office
# Theme" palette via the Design / ThemeColorScheme object model. The
# Colors() collection follows a:clrScheme's element order: dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink. RGB values use VBA's standard
# 0xBBGGRR (long) packing, matching what PowerPoint's RGB() builds.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$colors = $design.SlideMaster.Theme.ThemeColorScheme

function BGR($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (replaces the former Integral scheme).
$colors.Colors(1).RGB  = BGR 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = BGR 0x44 0x54 0x6A   # dk2
$colors.Colors(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2
$colors.Colors(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1
$colors.Colors(6).RGB  = BGR 0xED 0x7D 0x31   # accent2
$colors.Colors(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3
$colors.Colors(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4
$colors.Colors(9).RGB  = BGR 0x44 0x72 0xC4   # accent5
$colors.Colors(10).RGB = BGR 0x70 0xAD 0x47   # accent6
$colors.Colors(11).RGB = BGR 0x05 0x63 0xC1   # hlink
$colors.Colors(12).RGB = BGR 0x95 0x4F 0x72   # folHlink
